# The edit swaps the data of row 2 and row 3 (everything except the columns
# that already held identical values in both rows), and moves the
# "Antal substrat" / "Substrat-beskrivning" values (columns AN/AO) from row 3
# to row 2, since after the swap it is row 2 that carries that record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain numeric columns: straightforward value swap ---------------------
$ws.Range("A2").Value2  = 2355256
$ws.Range("A3").Value2  = 705191

$ws.Range("B2").Value2  = 76907
$ws.Range("B3").Value2  = 101201

$ws.Range("E2").Value2  = 6436
$ws.Range("E3").Value2  = 1910

$ws.Range("Q2").Value2  = 418748.8210844172
$ws.Range("Q3").Value2  = 418681.7318415901

$ws.Range("R2").Value2  = 6505592.199708293
$ws.Range("R3").Value2  = 6505423.5425952

$ws.Range("S2").Value2  = 10
$ws.Range("S3").Value2  = 100

# --- Plain text columns: straightforward value swap -------------------------
$ws.Range("D2").Value2  = "LC"
$ws.Range("D3").Value2  = "NT"

$ws.Range("F2").Value2  = "Gulpudrad spiklav"
$ws.Range("F3").Value2  = "Dvärglin"

$ws.Range("G2").Value2  = "Calicium adspersum"
$ws.Range("G3").Value2  = "Radiola linoides"

$ws.Range("H2").Value2  = "Pers."
$ws.Range("H3").Value2  = "Roth"

$ws.Range("P2").Value2  = "Årnäs, träd 13850 (ek), Vg"
$ws.Range("P3").Value2  = "Stranden NV om FUNKES, Vg"

$ws.Range("AW2").Value2 = "Lars Sjögren"
$ws.Range("AW3").Value2 = "Mora Aronsson"

$ws.Range("AX2").Value2 = "Andreas Furustam"
$ws.Range("AX3").Value2 = "Olof Janson"

$ws.Range("AY2").Value2 = "Epifyt-inv. Särskilt skyddsvärda träd 2005-2007"
$ws.Range("AY3").Value2 = "Västergötlands flora 2002"

# --- Date-like text columns ---------------------------------------------
# These hold plain text that LOOKS like a date ("yyyy-mm-dd"). Assigning such
# a string straight to .Value/.Value2 makes Excel re-interpret it as a real
# date serial number, which is not what the source data has (it is stored as
# literal text). Temporarily forcing a "Text" number format suppresses that
# auto-conversion; resetting the style back to Normal afterwards keeps the
# cell free of any left-over explicit formatting, matching the original file.
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value2 = "2005-08-16"
$ws.Range("Y2").Style = "Normal"

$ws.Range("Y3").NumberFormat = "@"
$ws.Range("Y3").Value2 = "1989-08-27"
$ws.Range("Y3").Style = "Normal"

$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value2 = "2005-08-16"
$ws.Range("AA2").Style = "Normal"

$ws.Range("AA3").NumberFormat = "@"
$ws.Range("AA3").Value2 = "1989-08-27"
$ws.Range("AA3").Style = "Normal"

# --- Columns AN/AO only exist (pre-edit) on row 3; move them to row 2 ------
$ws.Range("AN2").Value2 = 1
$ws.Range("AO2").Value2 = "1 substratenheter # ek"

$ws.Range("AN3").ClearContents()
$ws.Range("AO3").ClearContents()
